$targetList = @(
    'aws.s3',
    'aws.ses',
    'base',
    'csv',
    'desktop',
    'excel',
    'external',
    'image',
    'io',
    'jms',
    'json',
    'macro',
    'mail',
    'number',
    'pdf',
    'rdbms',
    'redis',
    'sms',
    'sound',
    'ssh',
    'step',
    'web',
    'webalert',
    'webcookie',
    'ws',
    'ws.async',
    'xml'
)

$webList = @(
    'assertAndClick(locator,label)',
    'assertAttribute(locator,attrName,value)',
    'assertAttributeContains(locator,attrName,contains)',
    'assertAttributeNotContains(locator,attrName,contains)',
    'assertAttributeNotPresent(locator,attrName)',
    'assertAttributePresent(locator,attrName)',
    'assertChecked(locator)',
    'assertContainCount(locator,text,count)',
    'assertCssNotPresent(locator,property)',
    'assertCssPresent(locator,property,value)',
    'assertElementByAttributes(nameValues)',
    'assertElementByText(locator,text)',
    'assertElementCount(locator,count)',
    'assertElementNotPresent(locator)',
    'assertElementPresent(locator)',
    'assertFocus(locator)',
    'assertFrameCount(count)',
    'assertFramePresent(frameName)',
    'assertIECompatMode()',
    'assertIENavtiveMode()',
    'assertLinkByLabel(label)',
    'assertNotChecked(locator)',
    'assertNotFocus(locator)',
    'assertNotText(locator,text)',
    'assertNotVisible(locator)',
    'assertOneMatch(locator)',
    'assertScrollbarHNotPresent(locator)',
    'assertScrollbarHPresent(locator)',
    'assertScrollbarVNotPresent(locator)',
    'assertScrollbarVPresent(locator)',
    'assertTable(locator,row,column,text)',
    'assertText(locator,text)',
    'assertTextContains(locator,text)',
    'assertTextCount(locator,text,count)',
    'assertTextList(locator,list,ignoreOrder)',
    'assertTextMatches(text,minMatch,scrollTo)',
    'assertTextNotPresent(text)',
    'assertTextOrder(locator,descending)',
    'assertTextPresent(text)',
    'assertTitle(text)',
    'assertValue(locator,value)',
    'assertValueOrder(locator,descending)',
    'assertVisible(locator)',
    'checkAll(locator)',
    'clearLocalStorage()',
    'click(locator)',
    'clickAndWait(locator,waitMs)',
    'clickByLabel(label)',
    'clickByLabelAndWait(label,waitMs)',
    'clickOffset(locator,x,y)',
    'clickWithKeys(locator,keys)',
    'close()',
    'closeAll()',
    'deselect(locator,text)',
    'deselectMulti(locator,array)',
    'dismissInvalidCert()',
    'dismissInvalidCertPopup()',
    'doubleClick(locator)',
    'doubleClickAndWait(locator,waitMs)',
    'doubleClickByLabel(label)',
    'doubleClickByLabelAndWait(label,waitMs)',
    'dragAndDrop(fromLocator,toLocator)',
    'dragTo(fromLocator,xOffset,yOffset)',
    'editLocalStorage(key,value)',
    'executeScript(var,script)',
    'focus(locator)',
    'goBack()',
    'goBackAndWait()',
    'maximizeWindow()',
    'mouseOver(locator)',
    'open(url)',
    'openAndWait(url,waitMs)',
    'openHttpBasic(url,username,password)',
    'openIgnoreTimeout(url)',
    'refresh()',
    'refreshAndWait()',
    'resizeWindow(width,height)',
    'saveAllWindowIds(var)',
    'saveAllWindowNames(var)',
    'saveAttribute(var,locator,attrName)',
    'saveAttributeList(var,locator,attrName)',
    'saveCount(var,locator)',
    'saveDivsAsCsv(headers,rows,cells,nextPage,file)',
    'saveElement(var,locator)',
    'saveElements(var,locator)',
    'saveLocalStorage(var,key)',
    'saveLocation(var)',
    'savePageAs(var,sessionIdName,url)',
    'savePageAsFile(sessionIdName,url,file)',
    'saveTableAsCsv(locator,nextPageLocator,file)',
    'saveText(var,locator)',
    'saveTextArray(var,locator)',
    'saveTextSubstringAfter(var,locator,delim)',
    'saveTextSubstringBefore(var,locator,delim)',
    'saveTextSubstringBetween(var,locator,start,end)',
    'saveValue(var,locator)',
    'scrollLeft(locator,pixel)',
    'scrollRight(locator,pixel)',
    'scrollTo(locator)',
    'select(locator,text)',
    'selectFrame(locator)',
    'selectMulti(locator,array)',
    'selectMultiOptions(locator)',
    'selectText(locator)',
    'selectWindow(winId)',
    'selectWindowAndWait(winId,waitMs)',
    'selectWindowByIndex(index)',
    'selectWindowByIndexAndWait(index,waitMs)',
    'toggleSelections(locator)',
    'type(locator,value)',
    'typeKeys(locator,value)',
    'uncheckAll(locator)',
    'unselectAllText()',
    'upload(fieldLocator,file)',
    'verifyContainText(locator,text)',
    'verifyText(locator,text)',
    'wait(waitMs)',
    'waitForElementPresent(locator)',
    'waitForPopUp(winId,waitMs)',
    'waitForTextPresent(text)',
    'waitForTitle(text)'
)

$externalList = @(
    'runJUnit(className)',
    'runProgram(programPathAndParams)',
    'runProgramNoWait(programPathAndParams)'
)

$macroList = @(
    'description()',
    'expects(var,default)',
    'produces(var,value)'
)


$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --------------------------------------------------------------------------
# 1) Insert a brand-new column at M ("macro"). This pushes the previous
#    M..AA columns (mail, number, pdf, rdbms, redis, sms, sound, ssh, step,
#    web, webalert, webcookie, ws, ws.async, xml) one column to the right,
#    becoming N..AB. Column-only insert (does not disturb A..L).
# --------------------------------------------------------------------------
$ws.Columns("M").Insert()

# --------------------------------------------------------------------------
# 2) Populate the new "macro" column (M) - header + 3 function names.
# --------------------------------------------------------------------------
$ws.Range("M1").Value = "macro"
for ($i = 0; $i -lt $macroList.Count; $i++) {
    $ws.Cells.Item($i + 2, 13).Value = $macroList[$i]
}

# --------------------------------------------------------------------------
# 3) "external" column (H) gains a 3rd function and fixes a typo in the
#    2nd ("programPathAndParms" -> "programPathAndParams"), plus adds the
#    brand new "runProgramNoWait(programPathAndParams)" entry.
# --------------------------------------------------------------------------
for ($i = 0; $i -lt $externalList.Count; $i++) {
    $ws.Cells.Item($i + 2, 8).Value = $externalList[$i]
}

# --------------------------------------------------------------------------
# 4) "target" column (A) - alphabetical list of all category/sheet names -
#    gains the new "macro" entry (inserted alphabetically between "json"
#    and "mail"). Rewrite the whole column from row 2 down with the final,
#    already-sorted list (now 27 entries instead of 26).
# --------------------------------------------------------------------------
for ($i = 0; $i -lt $targetList.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $targetList[$i]
}

# --------------------------------------------------------------------------
# 5) "web" column (now W, after the column insert above) - alphabetical
#    list of web commands - gains two new entries: "clickOffset(locator,x,y)"
#    and "saveAttributeList(var,locator,attrName)". Rewrite the whole
#    column from row 2 down with the final, already-sorted list (now 121
#    entries instead of 119).
# --------------------------------------------------------------------------
for ($i = 0; $i -lt $webList.Count; $i++) {
    $ws.Cells.Item($i + 2, 23).Value = $webList[$i]
}
